$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target keyword/search-volume table (rows 2-215); row 1 header and sheet structure unchanged.
$data = @(
    @('tiles', 11200),
    @('waterproof', 7480),
    @('waterproofing', 6890),
    @('terrace', 3900),
    @('bathroom', 3340),
    @('tile', 2060),
    @('tiler', 1900),
    @('waterproofer', 1900),
    @('water', 1420),
    @('proof', 1370),
    @('roof', 1360),
    @('floor', 1260),
    @('ceiling', 1050),
    @('shower', 1020),
    @('wall', 970),
    @('tiling', 950),
    @('tiled', 880),
    @('stick', 690),
    @('panels', 660),
    @('fixit', 540),
    @('dr', 540),
    @('gap', 530),
    @('peel', 520),
    @('vinyl', 470),
    @('weather', 460),
    @('on', 410),
    @('aquadefense', 390),
    @('price', 390),
    @('resistant', 330),
    @('without', 300),
    @('laticrete', 280),
    @('sheets', 280),
    @('removing', 270),
    @('myk', 260),
    @('over', 220),
    @('flooring', 220),
    @('tape', 200),
    @('after', 180),
    @('liquid', 170),
    @('walls', 170),
    @('best', 150),
    @('roff', 140),
    @('carpet', 140),
    @('proofing', 130),
    @('pvc', 120),
    @('basement', 120),
    @('bath', 110),
    @('fill', 110),
    @('kitchen', 100),
    @('fevimate', 90),
    @('backsplash', 90),
    @('kag', 80),
    @('suspended', 80),
    @('''false', 70),
    @('under', 70),
    @('existing', 70),
    @('interlocking', 70),
    @('palisade', 70),
    @('repellent', 60),
    @('ceramic', 60),
    @('of', 60),
    @('balcony', 60),
    @('tank', 50),
    @('tilesafe', 50),
    @('pidilite', 50),
    @('black', 50),
    @('ezee', 50),
    @('n', 50),
    @('kwik', 50),
    @('in', 50),
    @('top', 50),
    @('showers', 50),
    @('outdoor', 50),
    @('moisture', 50),
    @('reform', 40),
    @('joint', 40),
    @('wood', 40),
    @('deck', 40),
    @('chemicals', 30),
    @('kerakoll', 30),
    @('cement', 30),
    @('rubber', 30),
    @('sticky', 30),
    @('look', 30),
    @('foam', 30),
    @('smartcore', 30),
    @('duma', 30),
    @('like', 30),
    @('that', 30),
    @('transparent', 20),
    @('hard', 20),
    @('dunlop', 20),
    @('densshield', 20),
    @('pool', 20),
    @('swimming', 20),
    @('diy', 20),
    @('thermaldry', 20),
    @('curb', 20),
    @('acoustic', 20),
    @('mildewproof', 20),
    @('stone', 20),
    @('garage', 20),
    @('bathtub', 20),
    @('2x2', 20),
    @('looks', 20),
    @('wet', 10),
    @('weatherproof', 10),
    @('clay', 10),
    @('terracotta', 10),
    @('behind', 10),
    @('edging', 10),
    @('plastic', 10),
    @('mosaic', 10),
    @('marble', 10),
    @('sink', 10),
    @('&', 10),
    @('glover', 10),
    @('tony', 10),
    @('polyurethane', 10),
    @('shopee', 10),
    @('sikatite', 10),
    @('revela', 10),
    @('pre', 10),
    @('life', 10),
    @('hydroban', 10),
    @('glazed', 10),
    @('undertile', 10),
    @('rx3000', 10),
    @('dumawall+', 10),
    @('seam', 10),
    @('way', 10),
    @('aqua', 10),
    @('spectratile', 10),
    @('steam', 10),
    @('goboard', 10),
    @('veil', 10),
    @('maiden''s', 10),
    @('athena', 10),
    @('exterior', 10),
    @('porcelain', 10),
    @('instead', 10),
    @('armstrong', 10),
    @('spc', 10),
    @('leaking', 10),
    @('livelynine', 10),
    @('floating', 10),
    @('laminated', 10),
    @('spectra', 10),
    @('chatham', 10),
    @('moon', 10),
    @('half', 10),
    @('to', 10),
    @('mat', 10),
    @('600x600', 10),
    @('canyon', 10),
    @('desert', 10),
    @('groutable', 10),
    @('wedi', 10),
    @('12x24', 10),
    @('2x4', 10),
    @('home', 0),
    @('tarazzo', 0),
    @('wipeable', 0),
    @('slate', 0),
    @('bed', 0),
    @('mud', 0),
    @('together', 0),
    @('snap', 0),
    @('self', 0),
    @('linoleum', 0),
    @('groutless', 0),
    @('cork', 0),
    @('dry', 0),
    @('thermal', 0),
    @('technokolla', 0),
    @('hydraflex', 0),
    @('tec', 0),
    @('sumter', 0),
    @('redifix', 0),
    @('no', 0),
    @('gray', 0),
    @('whitfield', 0),
    @('msi', 0),
    @('jackoboard', 0),
    @('impervious', 0),
    @('hydrathane', 0),
    @('hybrid', 0),
    @('hanwood', 0),
    @('white', 0),
    @('pro', 0),
    @('smooth', 0),
    @('menards', 0),
    @('dukkaboard', 0),
    @('''3', 0),
    @('dampfix', 0),
    @('ampera', 0),
    @('bondera', 0),
    @('''13', 0),
    @('a108', 0),
    @('ansi', 0),
    @('absolute', 0),
    @('patio', 0),
    @('vv032', 0),
    @('coffer', 0),
    @('icon', 0),
    @('genesis', 0),
    @('concrete', 0),
    @('sedona', 0),
    @('coretec', 0),
    @('blonde', 0),
    @('bayhill', 0),
    @('everlife', 0),
    @('retile', 0),
    @('cost', 0)
)

$numRows = $data.Count
$values = New-Object 'object[,]' $numRows,2
for ($i = 0; $i -lt $numRows; $i++) {
    $values[$i, 0] = $data[$i][0]
    $values[$i, 1] = $data[$i][1]
}

$startRow = 2
$endRow = $startRow + $numRows - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 2))
$rng.Value = $values

# The table shrank by one row (216 -> 215); clear the now-unused trailing row.
$ws.Rows.Item(216).ClearContents()
